$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 417.53845
$ws.Cells.Item(12, 9).Value = 417.53845
$ws.Cells.Item(12, 11).Value = 417.53845
$ws.Cells.Item(12, 13).Value = -247.53845
$ws.Cells.Item(32, 8).Value = 13606.407
$ws.Cells.Item(32, 9).Value = 13969.4
$ws.Cells.Item(32, 11).Value = 13969.4
$ws.Cells.Item(32, 13).Value = -13643.4
$ws.Cells.Item(43, 8).Value = 1313
$ws.Cells.Item(43, 10).Value = 1059.4
$ws.Cells.Item(43, 12).Value = 1059.4
$ws.Cells.Item(43, 14).Value = -1197.4
$ws.Cells.Item(64, 8).Value = 444452770
$ws.Cells.Item(64, 9).Value = 444452770
$ws.Cells.Item(64, 11).Value = 444452770
$ws.Cells.Item(64, 13).Value = -444452522
$ws.Cells.Item(67, 8).Value = 444452770
$ws.Cells.Item(67, 9).Value = 444452770
$ws.Cells.Item(67, 11).Value = 444452770
$ws.Cells.Item(67, 13).Value = -444451912
$ws.Cells.Item(96, 8).Value = 58824976
$ws.Cells.Item(96, 9).Value = 1480.2
$ws.Cells.Item(96, 10).Value = 142858540
$ws.Cells.Item(96, 11).Value = 4440.6
$ws.Cells.Item(96, 12).Value = 428575620
$ws.Cells.Item(96, 13).Value = -3067.6
$ws.Cells.Item(96, 14).Value = -428578366
$ws.Cells.Item(100, 8).Value = 15383.571
$ws.Cells.Item(100, 10).Value = 23249.75
$ws.Cells.Item(100, 12).Value = 23249.75
$ws.Cells.Item(100, 14).Value = -24331.75
$ws.Cells.Item(118, 8).Value = 1067.5555
$ws.Cells.Item(118, 9).Value = 1067.5555
$ws.Cells.Item(118, 11).Value = 3202.6665
$ws.Cells.Item(118, 13).Value = -1545.6665
$ws.Cells.Item(129, 8).Value = 30304826
$ws.Cells.Item(129, 9).Value = 1615.1666
$ws.Cells.Item(129, 11).Value = 4845.4998
$ws.Cells.Item(129, 13).Value = 154.5002000000004
$ws.Cells.Item(132, 8).Value = 522816.25
$ws.Cells.Item(132, 9).Value = 924880.1
$ws.Cells.Item(132, 11).Value = 2774640.3
$ws.Cells.Item(132, 13).Value = -2772110.3
$ws.Cells.Item(138, 8).Value = 3666.948
$ws.Cells.Item(138, 9).Value = 1416.24
$ws.Cells.Item(138, 10).Value = 4459.4507
$ws.Cells.Item(138, 11).Value = 4248.72
$ws.Cells.Item(138, 12).Value = 13378.3521
$ws.Cells.Item(138, 13).Value = 891.2799999999997
$ws.Cells.Item(138, 14).Value = -23658.3521

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 10143.177
$ws.Cells.Item(2, 9).Value = 6036.4165
$ws.Cells.Item(2, 10).Value = 19999.4
$ws.Cells.Item(2, 11).Value = 6036.4165
$ws.Cells.Item(2, 12).Value = 19999.4
$ws.Cells.Item(2, 13).Value = -5923.4165
$ws.Cells.Item(2, 14).Value = -20225.4
$ws.Cells.Item(97, 8).Value = 517
$ws.Cells.Item(97, 9).Value = 425.73685
$ws.Cells.Item(97, 11).Value = 425.73685
$ws.Cells.Item(97, 13).Value = 70.26315
$ws.Cells.Item(116, 8).Value = 10143.177
$ws.Cells.Item(116, 9).Value = 6036.4165
$ws.Cells.Item(116, 10).Value = 19999.4
$ws.Cells.Item(116, 11).Value = 6036.4165
$ws.Cells.Item(116, 12).Value = 19999.4
$ws.Cells.Item(116, 13).Value = -3742.4165
$ws.Cells.Item(116, 14).Value = -24587.4
$ws.Cells.Item(132, 8).Value = 1346164.8
$ws.Cells.Item(132, 10).Value = 162333.17
$ws.Cells.Item(132, 12).Value = 486999.51
$ws.Cells.Item(132, 14).Value = -492059.51

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 10143.177
$ws.Cells.Item(3, 9).Value = 6036.4165
$ws.Cells.Item(3, 10).Value = 19999.4
$ws.Cells.Item(3, 11).Value = 6036.4165
$ws.Cells.Item(3, 12).Value = 19999.4
$ws.Cells.Item(3, 13).Value = -5922.4165
$ws.Cells.Item(3, 14).Value = -20227.4
$ws.Cells.Item(20, 8).Value = 8930628
$ws.Cells.Item(20, 9).Value = 15874795
$ws.Cells.Item(20, 10).Value = 2413.7144
$ws.Cells.Item(20, 11).Value = 15874795
$ws.Cells.Item(20, 12).Value = 2413.7144
$ws.Cells.Item(20, 13).Value = -15874548
$ws.Cells.Item(20, 14).Value = -2907.7144
$ws.Cells.Item(22, 8).Value = 2469.7144
$ws.Cells.Item(22, 9).Value = 3022.5
$ws.Cells.Item(22, 10).Value = 1732.6666
$ws.Cells.Item(22, 11).Value = 3022.5
$ws.Cells.Item(22, 12).Value = 1732.6666
$ws.Cells.Item(22, 13).Value = -2849.5
$ws.Cells.Item(22, 14).Value = -2078.6666
$ws.Cells.Item(86, 8).Value = 5805.2
$ws.Cells.Item(86, 10).Value = 7476.1665
$ws.Cells.Item(86, 12).Value = 7476.1665
$ws.Cells.Item(86, 14).Value = -9722.166499999999
$ws.Cells.Item(89, 8).Value = 5805.2
$ws.Cells.Item(89, 10).Value = 7476.1665
$ws.Cells.Item(89, 12).Value = 37380.8325
$ws.Cells.Item(89, 14).Value = -48612.8325
$ws.Cells.Item(94, 8).Value = 3940.8462
$ws.Cells.Item(94, 9).Value = 930
$ws.Cells.Item(94, 11).Value = 930
$ws.Cells.Item(94, 13).Value = -479
$ws.Cells.Item(135, 8).Value = 120780
$ws.Cells.Item(135, 10).Value = 120780
$ws.Cells.Item(135, 12).Value = 120780
$ws.Cells.Item(135, 14).Value = -130920

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 20003918
$ws.Cells.Item(16, 9).Value = 31252694
$ws.Cells.Item(16, 11).Value = 31252694
$ws.Cells.Item(16, 13).Value = -31252407
$ws.Cells.Item(105, 8).Value = 29415732
$ws.Cells.Item(105, 9).Value = 32261736
$ws.Cells.Item(105, 10).Value = 7033.3335
$ws.Cells.Item(105, 11).Value = 32261736
$ws.Cells.Item(105, 12).Value = 7033.3335
$ws.Cells.Item(105, 13).Value = -32259989
$ws.Cells.Item(105, 14).Value = -10527.3335
$ws.Cells.Item(113, 8).Value = 20003918
$ws.Cells.Item(113, 9).Value = 31252694
$ws.Cells.Item(113, 11).Value = 31252694
$ws.Cells.Item(113, 13).Value = -31250524
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 14).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 4074468
$ws.Cells.Item(4, 9).Value = 4703029
$ws.Cells.Item(4, 11).Value = 14109087
$ws.Cells.Item(4, 13).Value = -14108975
$ws.Cells.Item(17, 8).Value = 2655.25
$ws.Cells.Item(17, 10).Value = 1498.4445
$ws.Cells.Item(17, 12).Value = 4495.333500000001
$ws.Cells.Item(17, 14).Value = -4833.333500000001
$ws.Cells.Item(68, 8).Value = 240754.77
$ws.Cells.Item(68, 10).Value = 315578.12
$ws.Cells.Item(68, 12).Value = 946734.36
$ws.Cells.Item(68, 14).Value = -948356.36
$ws.Cells.Item(71, 8).Value = 240754.77
$ws.Cells.Item(71, 10).Value = 315578.12
$ws.Cells.Item(71, 12).Value = 2840203.08
$ws.Cells.Item(71, 14).Value = -2848315.08
$ws.Cells.Item(107, 8).Value = 4467.919
$ws.Cells.Item(107, 10).Value = 5389.276
$ws.Cells.Item(107, 12).Value = 16167.828
$ws.Cells.Item(107, 14).Value = -20007.828
$ws.Cells.Item(140, 8).Value = 98488740
$ws.Cells.Item(140, 9).Value = 98488740
$ws.Cells.Item(140, 11).Value = 295466220
$ws.Cells.Item(140, 13).Value = -295461040
$ws.Cells.Item(46, 8).Value = 3312.3076
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 3312.3076
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 9936.9228
$ws.Cells.Item(46, 13).ClearContents()
$ws.Cells.Item(46, 14).Value = -10118.9228

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 9747.25
$ws.Cells.Item(70, 9).Value = 9120.700000000001
$ws.Cells.Item(70, 10).Value = 10791.5
$ws.Cells.Item(70, 11).Value = 9120.700000000001
$ws.Cells.Item(70, 12).Value = 10791.5
$ws.Cells.Item(70, 13).Value = -8850.700000000001
$ws.Cells.Item(70, 14).Value = -11331.5
$ws.Cells.Item(73, 8).Value = 9747.25
$ws.Cells.Item(73, 9).Value = 9120.700000000001
$ws.Cells.Item(73, 10).Value = 10791.5
$ws.Cells.Item(73, 11).Value = 9120.700000000001
$ws.Cells.Item(73, 12).Value = 10791.5
$ws.Cells.Item(73, 13).Value = -8184.700000000001
$ws.Cells.Item(73, 14).Value = -12663.5
$ws.Cells.Item(97, 8).Value = 1090.5294
$ws.Cells.Item(97, 9).Value = 1059.24
$ws.Cells.Item(97, 11).Value = 1059.24
$ws.Cells.Item(97, 13).Value = -563.24
$ws.Cells.Item(122, 8).Value = 7753
$ws.Cells.Item(122, 9).Value = 7999.5
$ws.Cells.Item(122, 10).Value = 7136.75
$ws.Cells.Item(122, 11).Value = 23998.5
$ws.Cells.Item(122, 12).Value = 21410.25
$ws.Cells.Item(122, 13).Value = -21548.5
$ws.Cells.Item(122, 14).Value = -26310.25
$ws.Cells.Item(132, 8).Value = 142972670
$ws.Cells.Item(132, 9).Value = 200160590
$ws.Cells.Item(132, 10).Value = 2873.5
$ws.Cells.Item(132, 11).Value = 600481770
$ws.Cells.Item(132, 12).Value = 8620.5
$ws.Cells.Item(132, 13).Value = -600479240
$ws.Cells.Item(132, 14).Value = -13680.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 3260
$ws.Cells.Item(55, 10).Value = 5381.273
$ws.Cells.Item(55, 12).Value = 5381.273
$ws.Cells.Item(55, 14).Value = -5727.273
$ws.Cells.Item(100, 8).Value = 3357.0688
$ws.Cells.Item(100, 9).Value = 7395.2856
$ws.Cells.Item(100, 11).Value = 7395.2856
$ws.Cells.Item(100, 13).Value = -6854.2856

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 91285.14
$ws.Cells.Item(46, 10).Value = 91285.14
$ws.Cells.Item(46, 12).Value = 91285.14
$ws.Cells.Item(46, 14).Value = -91747.14
$ws.Cells.Item(62, 8).Value = 10957.417
$ws.Cells.Item(62, 9).Value = 10449.4
$ws.Cells.Item(62, 10).Value = 11320.286
$ws.Cells.Item(62, 11).Value = 10449.4
$ws.Cells.Item(62, 12).Value = 11320.286
$ws.Cells.Item(62, 13).Value = -9825.4
$ws.Cells.Item(62, 14).Value = -12568.286
$ws.Cells.Item(65, 8).Value = 10957.417
$ws.Cells.Item(65, 9).Value = 10449.4
$ws.Cells.Item(65, 10).Value = 11320.286
$ws.Cells.Item(65, 11).Value = 52247
$ws.Cells.Item(65, 12).Value = 56601.43
$ws.Cells.Item(65, 13).Value = -49127
$ws.Cells.Item(65, 14).Value = -62841.43
$ws.Cells.Item(126, 8).Value = 5798.25
$ws.Cells.Item(126, 9).Value = 3524.2104
$ws.Cells.Item(126, 10).Value = 14439.6
$ws.Cells.Item(126, 11).Value = 10572.6312
$ws.Cells.Item(126, 12).Value = 43318.8
$ws.Cells.Item(126, 13).Value = -8102.6312
$ws.Cells.Item(126, 14).Value = -48258.8
$ws.Cells.Item(134, 8).Value = 91285.14
$ws.Cells.Item(134, 10).Value = 91285.14
$ws.Cells.Item(134, 12).Value = 273855.42
$ws.Cells.Item(134, 14).Value = -278925.42
